$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 55557190
$ws.Range("I100").Value = 83334290
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 83334290
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -83333749
$ws.Range("N100").Value = -4082

$ws.Range("H107").Value = 11366788
$ws.Range("I107").Value = 12500942
$ws.Range("J107").Value = 25250
$ws.Range("K107").Value = 12500942
$ws.Range("L107").Value = 25250
$ws.Range("M107").Value = -12499022
$ws.Range("N107").Value = -29090

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4168.202
$ws.Range("I32").Value = 3549.5823
$ws.Range("J32").Value = 7426.2666
$ws.Range("K32").Value = 3549.5823
$ws.Range("L32").Value = 7426.2666
$ws.Range("M32").Value = -3262.5823
$ws.Range("N32").Value = -8000.2666

$ws.Range("H61").Value = 55000
$ws.Range("I61").Value = 10000
$ws.Range("J61").Value = 100000
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 100000
$ws.Range("M61").Value = -9788
$ws.Range("N61").Value = -100424

$ws.Range("H74").Value = 1128.5714
$ws.Range("I74").Value = 856.7442
$ws.Range("J74").Value = 2027.6923
$ws.Range("K74").Value = 856.7442
$ws.Range("L74").Value = 2027.6923
$ws.Range("M74").Value = 17.25580000000002
$ws.Range("N74").Value = -3775.6923

$ws.Range("H77").Value = 1128.5714
$ws.Range("I77").Value = 856.7442
$ws.Range("J77").Value = 2027.6923
$ws.Range("K77").Value = 4283.721
$ws.Range("L77").Value = 10138.4615
$ws.Range("M77").Value = 84.27900000000045
$ws.Range("N77").Value = -18874.4615

$ws.Range("H110").Value = 1150.3334
$ws.Range("I110").Value = 725.5
$ws.Range("K110").Value = 725.5
$ws.Range("M110").Value = 1319.5

$ws.Range("H132").Value = 7632.1875
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 7474.3335
$ws.Range("K132").Value = 30000
$ws.Range("L132").Value = 22423.0005
$ws.Range("M132").Value = -27470
$ws.Range("N132").Value = -27483.0005

$ws.Range("H136").Value = 55000
$ws.Range("I136").Value = 10000
$ws.Range("J136").Value = 100000
$ws.Range("K136").Value = 30000
$ws.Range("L136").Value = 300000
$ws.Range("M136").Value = -27450
$ws.Range("N136").Value = -305100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 561
$ws.Range("I5").Value = 561
$ws.Range("K5").Value = 561
$ws.Range("M5").Value = -448

$ws.Range("H86").Value = 2690.0908
$ws.Range("I86").Value = 2461.375
$ws.Range("J86").Value = 3300
$ws.Range("K86").Value = 2461.375
$ws.Range("L86").Value = 3300
$ws.Range("M86").Value = -1338.375
$ws.Range("N86").Value = -5546

$ws.Range("H89").Value = 2690.0908
$ws.Range("I89").Value = 2461.375
$ws.Range("J89").Value = 3300
$ws.Range("K89").Value = 12306.875
$ws.Range("L89").Value = 16500
$ws.Range("M89").Value = -6690.875
$ws.Range("N89").Value = -27732

$ws.Range("H107").Value = 1600
$ws.Range("I107").Value = 1600
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1600
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 320
$ws.Range("N107").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 76923070
$ws.Range("I16").Value = 76923070
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 76923070
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -76922783
$ws.Range("N16").Value = $null

$ws.Range("H31").Value = 3882.88
$ws.Range("I31").Value = 1616.9615
$ws.Range("J31").Value = 5085.204
$ws.Range("K31").Value = 1616.9615
$ws.Range("L31").Value = 5085.204
$ws.Range("M31").Value = -1321.9615
$ws.Range("N31").Value = -5675.204

$ws.Range("H34").Value = 3882.88
$ws.Range("I34").Value = 1616.9615
$ws.Range("J34").Value = 5085.204
$ws.Range("K34").Value = 1616.9615
$ws.Range("L34").Value = 5085.204
$ws.Range("M34").Value = -1414.9615
$ws.Range("N34").Value = -5489.204

$ws.Range("H99").Value = 9629770
$ws.Range("I99").Value = 18420
$ws.Range("J99").Value = 15636864
$ws.Range("K99").Value = 18420
$ws.Range("L99").Value = 15636864
$ws.Range("M99").Value = -16922
$ws.Range("N99").Value = -15639860

$ws.Range("H107").Value = 673
$ws.Range("I107").Value = 610.1667
$ws.Range("J107").Value = 731
$ws.Range("K107").Value = 610.1667
$ws.Range("L107").Value = 731
$ws.Range("M107").Value = 1309.8333
$ws.Range("N107").Value = -4571

$ws.Range("H113").Value = 76923070
$ws.Range("I113").Value = 76923070
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 76923070
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -76920900
$ws.Range("N113").Value = $null

$ws.Range("H122").Value = 837.1667
$ws.Range("I122").Value = 837.1667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2511.5001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -61.5001000000002
$ws.Range("N122").Value = $null

$ws.Range("H126").Value = 9629770
$ws.Range("I126").Value = 18420
$ws.Range("J126").Value = 15636864
$ws.Range("K126").Value = 55260
$ws.Range("L126").Value = 46910592
$ws.Range("M126").Value = -52790
$ws.Range("N126").Value = -46915532

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").Value = $null

$ws.Range("H37").Value = 42375
$ws.Range("J37").Value = 42375
$ws.Range("L37").Value = 127125
$ws.Range("N37").Value = -127349

$ws.Range("H68").Value = 2898.5557
$ws.Range("I68").Value = 4094.9033
$ws.Range("J68").Value = 1739.5938
$ws.Range("K68").Value = 12284.7099
$ws.Range("L68").Value = 5218.7814
$ws.Range("M68").Value = -11473.7099
$ws.Range("N68").Value = -6840.7814

$ws.Range("H71").Value = 2898.5557
$ws.Range("I71").Value = 4094.9033
$ws.Range("J71").Value = 1739.5938
$ws.Range("K71").Value = 36854.1297
$ws.Range("L71").Value = 15656.3442
$ws.Range("M71").Value = -32798.1297
$ws.Range("N71").Value = -23768.3442

$ws.Range("H107").Value = 840.80646
$ws.Range("I107").Value = 290.0909
$ws.Range("J107").Value = 1143.7
$ws.Range("K107").Value = 870.2727
$ws.Range("L107").Value = 3431.1
$ws.Range("M107").Value = 1049.7273
$ws.Range("N107").Value = -7271.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 5500000
$ws.Range("I24").Value = 5500000
$ws.Range("K24").Value = 5500000
$ws.Range("M24").Value = -5499827

$ws.Range("H80").Value = 2669.423
$ws.Range("I80").Value = 2619.0625
$ws.Range("J80").Value = 2750
$ws.Range("K80").Value = 2619.0625
$ws.Range("L80").Value = 2750
$ws.Range("M80").Value = -1621.0625
$ws.Range("N80").Value = -4746

$ws.Range("H83").Value = 2669.423
$ws.Range("I83").Value = 2619.0625
$ws.Range("J83").Value = 2750
$ws.Range("K83").Value = 13095.3125
$ws.Range("L83").Value = 13750
$ws.Range("M83").Value = -8103.3125
$ws.Range("N83").Value = -23734

$ws.Range("H107").Value = 4189.1113
$ws.Range("I107").Value = 10000
$ws.Range("J107").Value = 3462.75
$ws.Range("K107").Value = 10000
$ws.Range("L107").Value = 3462.75
$ws.Range("M107").Value = -8080
$ws.Range("N107").Value = -7302.75

$ws.Range("H132").Value = 5596.6
$ws.Range("I132").Value = 11666.667
$ws.Range("J132").Value = 4079.0833
$ws.Range("K132").Value = 35000.001
$ws.Range("L132").Value = 12237.2499
$ws.Range("M132").Value = -32470.001
$ws.Range("N132").Value = -17297.2499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 495
$ws.Range("I16").Value = 495
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 495
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -325
$ws.Range("N16").Value = $null

$ws.Range("H61").Value = 2150
$ws.Range("I61").Value = 2150
$ws.Range("K61").Value = 2150
$ws.Range("M61").Value = -1948

$ws.Range("H93").Value = 71459030
$ws.Range("I93").Value = 41640.2
$ws.Range("K93").Value = 41640.2
$ws.Range("M93").Value = -40392.2

$ws.Range("H113").Value = 2150
$ws.Range("I113").Value = 2150
$ws.Range("K113").Value = 2150
$ws.Range("M113").Value = 20

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 125000536
$ws.Range("I107").Value = 250000480
$ws.Range("J107").Value = 600.25
$ws.Range("K107").Value = 750001440
$ws.Range("L107").Value = 1800.75
$ws.Range("M107").Value = -749999520
$ws.Range("N107").Value = -5640.75

$ws.Range("H113").Value = 2209
$ws.Range("I113").Value = 1507
$ws.Range("J113").Value = 3437.5
$ws.Range("K113").Value = 4521
$ws.Range("L113").Value = 10312.5
$ws.Range("M113").Value = -2351
$ws.Range("N113").Value = -14652.5

$ws.Range("H132").Value = 2708.0476
$ws.Range("I132").Value = 3027.5454
$ws.Range("J132").Value = 2356.6
$ws.Range("K132").Value = 9082.636200000001
$ws.Range("L132").Value = 7069.799999999999
$ws.Range("M132").Value = -6552.636200000001
$ws.Range("N132").Value = -12129.8
